# Auto-generated Excel COM-interop script
# Applies numeric cell updates across multiple worksheets as described by the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 99.75
$ws.Range("I6").Value = 99.75
$ws.Range("K6").Value = 299.25
$ws.Range("M6").Value = -187.25
$ws.Range("H74").Value = 166676670
$ws.Range("I74").Value = 250005000
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 250005000
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -250004064
$ws.Range("N74").Value = -21872
$ws.Range("H77").Value = 166676670
$ws.Range("I77").Value = 250005000
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 1250025000
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -1250020320
$ws.Range("N77").Value = -109360
$ws.Range("H137").Value = 2447.3914
$ws.Range("I137").Value = 2524.4285
$ws.Range("J137").Value = 2327.5557
$ws.Range("K137").Value = 7573.2855
$ws.Range("L137").Value = 6982.6671
$ws.Range("M137").Value = -5023.2855
$ws.Range("N137").Value = -12082.6671
$ws.Range("H138").Value = 3288.3103
$ws.Range("I138").Value = 922.2857
$ws.Range("J138").Value = 5496.6
$ws.Range("K138").Value = 2766.8571
$ws.Range("L138").Value = 16489.8
$ws.Range("M138").Value = 2373.1429
$ws.Range("N138").Value = -26769.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1457270.9
$ws.Range("I32").Value = 1509093.8
$ws.Range("K32").Value = 1509093.8
$ws.Range("M32").Value = -1508806.8
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H110").Value = 12821307
$ws.Range("I110").Value = 810.913
$ws.Range("K110").Value = 810.913
$ws.Range("M110").Value = 1234.087
$ws.Range("H122").Value = 26811.334
$ws.Range("I122").Value = 34883.668
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 104651.004
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -102201.004
$ws.Range("N122").Value = -36900.001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H99").Value = 3955368
$ws.Range("I99").Value = 2673.4119
$ws.Range("J99").Value = 15154670
$ws.Range("K99").Value = 2673.4119
$ws.Range("L99").Value = 15154670
$ws.Range("M99").Value = -1175.4119
$ws.Range("N99").Value = -15157666
$ws.Range("H107").Value = 43272744
$ws.Range("I107").Value = 51138856
$ws.Range("K107").Value = 51138856
$ws.Range("M107").Value = -51136936
$ws.Range("H134").Value = 6022.927
$ws.Range("I134").Value = 2278.0715
$ws.Range("J134").Value = 7964.7036
$ws.Range("K134").Value = 6834.2145
$ws.Range("L134").Value = 23894.1108
$ws.Range("M134").Value = -4299.2145
$ws.Range("N134").Value = -28964.1108

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 112990.336
$ws.Range("J122").Value = 336004.66
$ws.Range("L122").Value = 1008013.98
$ws.Range("N122").Value = -1012913.98

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 13889196
$ws.Range("I14").Value = 13889196
$ws.Range("K14").Value = 41667588
$ws.Range("M14").Value = -41667415
$ws.Range("H98").Value = 3521.3845
$ws.Range("J98").Value = 3755.111
$ws.Range("L98").Value = 11265.333
$ws.Range("N98").Value = -14261.333
$ws.Range("H107").Value = 14286055
$ws.Range("J107").Value = 18182142
$ws.Range("L107").Value = 54546426
$ws.Range("N107").Value = -54550266
$ws.Range("H114").Value = 854.25
$ws.Range("I114").Value = 782.5
$ws.Range("J114").Value = 897.3
$ws.Range("K114").Value = 2347.5
$ws.Range("L114").Value = 2691.9
$ws.Range("M114").Value = 906.5
$ws.Range("N114").Value = -9199.9
$ws.Range("H121").Value = 6250989
$ws.Range("I121").Value = 65
$ws.Range("J121").Value = 8334630.5
$ws.Range("K121").Value = 195
$ws.Range("L121").Value = 25003891.5
$ws.Range("M121").Value = 1115
$ws.Range("N121").Value = -25006511.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H132").Value = 5724.522
$ws.Range("I132").Value = 2346.4285
$ws.Range("K132").Value = 7039.2855
$ws.Range("M132").Value = -4509.2855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4200.2
$ws.Range("I2").Value = 1001
$ws.Range("K2").Value = 1001
$ws.Range("M2").Value = -889
$ws.Range("H22").Value = 1555.1936
$ws.Range("I22").Value = 1116
$ws.Range("J22").Value = 3385.1667
$ws.Range("K22").Value = 1116
$ws.Range("L22").Value = 3385.1667
$ws.Range("M22").Value = -821
$ws.Range("N22").Value = -3975.1667
$ws.Range("H27").Value = 1555.1936
$ws.Range("I27").Value = 1116
$ws.Range("J27").Value = 3385.1667
$ws.Range("K27").Value = 1116
$ws.Range("L27").Value = 3385.1667
$ws.Range("M27").Value = -1009
$ws.Range("N27").Value = -3599.1667
$ws.Range("H46").Value = 2251.5
$ws.Range("J46").Value = 2413.3333
$ws.Range("L46").Value = 2413.3333
$ws.Range("N46").Value = -2789.3333
$ws.Range("H55").Value = 27778040
$ws.Range("I55").Value = 58823596
$ws.Range("J55").Value = 434.78946
$ws.Range("K55").Value = 58823596
$ws.Range("L55").Value = 434.78946
$ws.Range("M55").Value = -58823423
$ws.Range("N55").Value = -780.78946
$ws.Range("H88").Value = 34750
$ws.Range("J88").Value = 34750
$ws.Range("L88").Value = 34750
$ws.Range("N88").Value = -35606
$ws.Range("H91").Value = 34750
$ws.Range("J91").Value = 34750
$ws.Range("L91").Value = 34750
$ws.Range("N91").Value = -37714
$ws.Range("H100").Value = 4582.6924
$ws.Range("I100").Value = 3508.2222
$ws.Range("K100").Value = 3508.2222
$ws.Range("M100").Value = -2967.2222

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 27059402
$ws.Range("I136").Value = 58824660
$ws.Range("K136").Value = 176473980
$ws.Range("M136").Value = -176471430

Write-Host "Applied all profit-sheet updates."
